$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=8.889220333333332; H=26.667661; I=0.5709282459015017; J=0.5709282459015017;
            K=3; M=0.2524303333333334; N=0.757291; O=0.0004152138066912919; P=0.0004153435393530768;
            Q=2.243908851816778; R=20.195179666351; S=0.0002370572903283445; T=0.0002371313583693735 }
    3  = @{ E=3; G=8.889220333333332; H=26.667661; I=0.5709282459015017; J=0.5709282459015017;
            K=3; M=0.3565143333333333; N=1.069543; O=0.0005864179297654724; P=0.0005866011547876679;
            Q=3.169134460991444; R=28.522210148923; S=0.0003348025600061912; T=0.0003349071683467185 }
    4  = @{ E=3; G=8.889220333333332; H=26.667661; I=0.5709282459015017; J=0.5709282459015017;
            K=3; M=343.6225723333334; N=1030.867717; O=0.5652127240000627; P=0.5653893235012776;
            Q=3054.536756977771; R=27490.83081279994; S=0.3226959090745654; T=0.3227967347180211 }
    5  = @{ E=3; G=8.889220333333332; H=26.667661; I=0.5709282459015017; J=0.5709282459015017;
            K=3; M=263.1514486666667; N=789.454346; O=0.432848593488691; P=0.4329838361017211;
            Q=2339.211208233856; R=21052.9008741047; S=0.2471254882214305; T=0.247202702049259 }
    6  = @{ E=3; G=8.889220333333332; H=26.667661; I=0.5709282459015017; J=0.5709282459015017;
            K=2; M=0.5696825; N=1.139365; O=0.000937050774789659; P=0.0006248957028606154;
            Q=5.064033262544166; R=30.384199575265; S=0.0005349887551713031; T=0.0003567706075055972 }
    7  = @{ E=3; G=6.680547666666667; H=20.041643; I=0.4290717540984982; J=0.4290717540984982;
            K=3; M=0.2524303333333334; N=0.757291; O=0.0004152138066912919; P=0.0004153435393530768;
            Q=1.686372874345889; R=15.177355869113; S=0.0001781565163629474; T=0.0001782121809837033 }
    8  = @{ E=3; G=6.680547666666667; H=20.041643; I=0.4290717540984982; J=0.4290717540984982;
            K=3; M=0.3565143333333333; N=1.069543; O=0.0005864179297654724; P=0.0005866011547876679;
            Q=2.381710997683222; R=21.435398979149; S=0.0002516153697592812; T=0.0002516939864409493 }
    9  = @{ E=3; G=6.680547666666667; H=20.041643; I=0.4290717540984982; J=0.4290717540984982;
            K=3; M=343.6225723333334; N=1030.867717; O=0.5652127240000627; P=0.5653893235012776;
            Q=2295.586973815448; R=20660.28276433903; S=0.2425168149254973; T=0.2425925887832564 }
    10 = @{ E=3; G=6.680547666666667; H=20.041643; I=0.4290717540984982; J=0.4290717540984982;
            K=3; M=263.1514486666667; N=789.454346; O=0.432848593488691; P=0.4329838361017211;
            Q=1757.995796370053; R=15821.96216733048; S=0.1857231052672604; T=0.1857811340524622 }
    11 = @{ E=3; G=6.680547666666667; H=20.041643; I=0.4290717540984982; J=0.4290717540984982;
            K=2; M=0.5696825; N=1.139365; O=0.000937050774789659; P=0.0006248957028606154;
            Q=3.805791096115833; R=22.834746576695; S=0.0004020620196183558; T=0.0002681250953550182 }
}

foreach ($rowNum in $data.Keys) {
    $rowData = $data[$rowNum]
    foreach ($col in $rowData.Keys) {
        $addr = "$col$rowNum"
        $ws.Range($addr).Value = $rowData[$col]
    }
}
